# "Generate Report for Handback" — mark the two handed-off files (the
# 2a8a0528... and d777324f... markdown files) as handed back / in sync,
# for both the zh-cn and de-de localization targets, and reflect that on
# the Overview sheet.

$wb = $excel.ActiveWorkbook

$urlMd1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/38a36c9885fb1af62d0228f32eb5ba585fcbeae8/e2e/2a8a0528-1b74-4041-a344-c4b731c644b3.md"
$urlMd2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/38a36c9885fb1af62d0228f32eb5ba585fcbeae8/e2e/d777324f-c81b-4580-bd8e-123f1e752530.md"

$hyperlinkColor = 15570276  # BGR for RGB(100,149,237) / #6495ED, matching the workbook's existing "HyperLink" cell style

# ---------------------------------------------------------------------
# Overview sheet: status moves from "Ready for handoff" to "Handed back"
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# Per-locale report sheets: zh-cn and de-de each get, for both rows
# (the two source files), a "Latest Target File" (I) and
# "Latest Handback File" (J) populated, plus a fresh
# "Latest Handback DateTime" (K).
# ---------------------------------------------------------------------
$locales = @(
    @{ Name = "zh-cn"; Suffix = "zh-cn"; HandbackTime = "2016-08-25 22:59:47" },
    @{ Name = "de-de"; Suffix = "de-de"; HandbackTime = "2016-08-25 22:59:54" }
)

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.Name)

    # Widen columns C (Status), I (Latest Target File), J (Latest Handback File)
    $ws.Columns.Item(3).ColumnWidth = 29.166666666666668
    $ws.Columns.Item(9).ColumnWidth = 39.166666666666664
    $ws.Columns.Item(10).ColumnWidth = 39.166666666666664

    # --- Row 2: 2a8a0528-1b74-4041-a344-c4b731c644b3.md ---
    $ws.Hyperlinks.Add($ws.Range("I2"), $urlMd1, "", "", "2a8a0528-1b74-4041-a344-c4b731c644b3.md") | Out-Null
    $ws.Range("I2").Font.Underline = $true
    $ws.Range("I2").Font.Color = $hyperlinkColor
    $ws.Range("J2").Value = "2a8a0528-1b74-4041-a344-c4b731c644b3.324e2b0f1554a297d6469b0ac684e22365e6686b." + $loc.Suffix + ".xlf"
    $ws.Range("K2").Value = $loc.HandbackTime

    # --- Row 3: d777324f-c81b-4580-bd8e-123f1e752530.md ---
    $ws.Hyperlinks.Add($ws.Range("I3"), $urlMd2, "", "", "d777324f-c81b-4580-bd8e-123f1e752530.md") | Out-Null
    $ws.Range("I3").Font.Underline = $true
    $ws.Range("I3").Font.Color = $hyperlinkColor
    $ws.Range("J3").Value = "d777324f-c81b-4580-bd8e-123f1e752530.1bfe4cb5a8237ead60d46a681c598765c5bbd202." + $loc.Suffix + ".xlf"
    $ws.Range("K3").Value = $loc.HandbackTime
}
